$d = $word.ActiveDocument

function Get-ParagraphEnd($pos) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Start -le $pos -and $p.Range.End -ge $pos) {
            return $p.Range.End
        }
    }
    return $pos
}

# --- Edit 1: split "{{NUMERO_PARECER}}" run into "{{" + "NUMERO_PARECER}}" ---
$r1 = $d.Content
$found1 = $r1.Find.Execute("{{NUMERO_PARECER}}", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) { throw "Could not find {{NUMERO_PARECER}} target" }
$pEnd1 = Get-ParagraphEnd($r1.Start)
$span1 = $d.Range($r1.Start, $pEnd1)
$xml1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t>{{</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t>NUMERO_PARECER}}</w:t></w:r><w:r w:rsidRPr="00EE2C2C"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t xml:space="preserve"> DA COMISSÃO DE EDUCAÇÃO, SAÚDE E ASSIST. SOCIAL </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$span1.InsertXML($xml1)

# --- Edit 2: split "{{NUMERO_PROJETO}}," run into "{{NUMERO_PROJETO}}" + "," ---
$r2 = $d.Content
$found2 = $r2.Find.Execute("{{NUMERO_PROJETO}},", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) { throw "Could not find {{NUMERO_PROJETO}}, target" }
$pEnd2 = Get-ParagraphEnd($r2.Start)
$span2 = $d.Range($r2.Start, $pEnd2)
$xml2 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t>{{NUMERO_PROJETO}}</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t>,</w:t></w:r><w:r w:rsidRPr="00EE2C2C"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t xml:space="preserve"> DE {{DATA_PROJETO}}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$span2.InsertXML($xml2)

# --- Edit 3: split "{{DATA_PROTOCOLO}}{{TEXTO_APRESENTACAO}}" run into two runs ---
$r3 = $d.Content
$found3 = $r3.Find.Execute("{{DATA_PROTOCOLO}}{{TEXTO_APRESENTACAO}}", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found3) { throw "Could not find {{DATA_PROTOCOLO}}{{TEXTO_APRESENTACAO}} target" }
$pEnd3 = Get-ParagraphEnd($r3.Start)
$span3 = $d.Range($r3.Start, $pEnd3)
$xml3 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>{{DATA_PROTOCOLO}}</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>{{TEXTO_APRESENTACAO}}</w:t></w:r><w:r w:rsidRPr="00583A72"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>. A proposição foi encaminhada à Comissão de Justiça e Redação, onde recebeu parecer favorável à sua tramitação. Em seguida, o Projeto foi encaminhado à Comissão de Obras, Serviço Públicos e Atividades Privadas, e, por fim, a esta Comissão de Educação, Saúde e Assistência Social, a fim de ser analisada, conforme previsto no artigo 47 do Regimento Interno.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$span3.InsertXML($xml3)

Write-Host "All edits applied"
